$wb = $excel.ActiveWorkbook

# Rename the two "prod" sheets, swapping their display names:
#  - the sheet currently named "prodTestSheetName" becomes "xprodTestSheetName"
#  - the sheet currently named "xxprodTestSheetName" becomes "prodTestSheetName"
$wb.Worksheets.Item("prodTestSheetName").Name = "xprodTestSheetName"
$wb.Worksheets.Item("xxprodTestSheetName").Name = "prodTestSheetName"

# Make the (newly renamed) "prodTestSheetName" sheet the active tab
$wb.Worksheets.Item("prodTestSheetName").Activate()
